# TC16_INS_CancerType-PancreasCancer.xlsx
# Commit: "automation API poc changes & INS complete regression suite 23 scripts"
#
# The "Program" query cell (B2) was rewritten: the "Website" column of the
# SQL query was changed from a straight `prg.website` projection to a CASE
# expression that prefers `prg.program_acronym` (falling back to
# `prg.program_link`) whenever `prg.program_link` is present. The selected
# cell in the sheet also moved from B5 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newProgramQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Pancreas Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

# Trim the single trailing newline that the PowerShell here-string picks up
# (the source cell text ends at the closing ';', with no trailing blank line).
$newProgramQuery = $newProgramQuery.TrimEnd("`r", "`n")

# Re-type the Program-tab query cell with the updated SQL. Re-applying the
# font size mirrors the (re-)formatting Excel recorded for this cell when it
# was edited in the UI.
$ws.Range("B2").Value = $newProgramQuery
$ws.Range("B2").Font.Size = 12

# Match the author's final selection/scroll position (was B5 with the view
# scrolled so row 5 was at the top; now C3, fully in view so no custom
# top-left is needed).
$ws.Range("C3").Select()

Write-Output "Updated Program query in B2 and moved selection to C3"
